# Leetcode_daily workbook update — add 6 new daily-log rows (311-316) that
# were previously blank placeholder rows, and nudge the window/pane view
# like Excel does after scrolling/selecting a new active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 311 - 583. Delete Operation for Two Strings (medium)
# ---------------------------------------------------------------------
$ws.Range("A311").Value = "583. Delete Operation for Two Strings"
$ws.Range("B311").Value = "medium"
$ws.Range("D311").Value = "DFS+记忆化: 找时间再仔细把DP写法学一下，看看三叶姐姐有没有专题 #难题呀字符串问题，需要下功夫"
$ws.Range("C310").Copy()
$ws.Range("C311").PasteSpecial(-4122)
$ws.Range("C311").Value = 44464

# ---------------------------------------------------------------------
# Row 312 - 513. Find Bottom Left Tree Value (medium)
# ---------------------------------------------------------------------
$ws.Range("A312").Value = "513. Find Bottom Left Tree Value`n"
$ws.Range("B312").Value = "medium"
$ws.Range("D312").Value = "BFS：板子送分题 DFS：使用先序遍历 维护nonlocal ans 和nonlocal maxdepth 当深度大于max就更新ans一次。"
$ws.Range("C310").Copy()
$ws.Range("C312").PasteSpecial(-4122)
$ws.Range("C312").Value = 44464

# ---------------------------------------------------------------------
# Row 313 - 371. Sum of Two Integers (medium)
# ---------------------------------------------------------------------
$ws.Range("A313").Value = "371. Sum of Two Integers`n"
$ws.Range("B313").Value = "medium"
$ws.Range("D313").Value = "bits manipulation： 与运算1与1为1 或运算1或0为1 异或运算相异为1 原理是利用异或运算为没有进位的加法，而与运算再<<1以后为进行，一直循环直到没有进位退出即可。"
$ws.Range("C310").Copy()
$ws.Range("C313").PasteSpecial(-4122)
$ws.Range("C313").Value = 44465

# ---------------------------------------------------------------------
# Row 314 - 297. Serialize and Deserialize Binary Tree (hard)
# ---------------------------------------------------------------------
$ws.Range("A314").Value = "297. Serialize and Deserialize Binary Tree`n"
$ws.Range("B314").Value = "hard"
$ws.Range("D314").Value = "DFS：序列化：直接前序遍历二叉树将其val生成字符串 运用字符串可拼接性质 反序列化：先split(,)切割为数组，然后每次弹出最前元素进行转化，如果为None就返回None，不是的话就做个root，然后递归找它的子节点。"
$ws.Range("C310").Copy()
$ws.Range("C314").PasteSpecial(-4122)
$ws.Range("C314").Value = 44465

# ---------------------------------------------------------------------
# Row 315 - 987. Vertical Order Traversal of a Binary Tree (hard)
# ---------------------------------------------------------------------
$ws.Range("A315").Value = "987. Vertical Order Traversal of a Binary Tree"
$ws.Range("B315").Value = "hard"
$ws.Range("D315").Value = "DFS+sort： dfs遍历, 得到col,row,value三元组 ,col 为第一关键字升序,row为第二关键字升序,value 为第三关键字升序, 同列存到字典,key为col,value为[val]"
$ws.Range("D184").Copy()
$ws.Range("C315").PasteSpecial(-4122)
$ws.Range("C315").Value = 44466

# ---------------------------------------------------------------------
# Row 316 - 639. Decode Ways II (hard)
# ---------------------------------------------------------------------
$ws.Range("A316").Value = "639. Decode Ways II"
$ws.Range("B316").Value = "hard"
$ws.Range("D316").Value = "划水过"
$ws.Range("D184").Copy()
$ws.Range("C316").PasteSpecial(-4122)
$ws.Range("C316").Value = 44466

# ---------------------------------------------------------------------
# Refresh the view: active cell / selection moved to D316 as the user
# scrolled the frozen pane down to keep row 313 at the top.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 313
$ws.Range("D316").Select()
